# Update "想去人数" (F column) values on the 展览, 演出 and 全部类型 sheets
# to reflect the latest scrape (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 283
$wsExhibit.Range("F3").Value = 193
$wsExhibit.Range("F4").Value = 2249
$wsExhibit.Range("F5").Value = 1712
$wsExhibit.Range("F6").Value = 326
$wsExhibit.Range("F8").Value = 780
$wsExhibit.Range("F9").Value = 162

# Sheet "演出" (performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 23

# Sheet "全部类型" (all types, combined list)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 283
$wsAll.Range("F3").Value = 193
$wsAll.Range("F4").Value = 2249
$wsAll.Range("F5").Value = 1712
$wsAll.Range("F6").Value = 326
$wsAll.Range("F7").Value = 23
$wsAll.Range("F9").Value = 780
$wsAll.Range("F10").Value = 162
